# Generate Report for Handoff
#
# A new handoff run produced a fresh report: the handed-off markdown file
# and its translated .xlf targets were regenerated under a new GUID/hash,
# and the handoff timestamps advanced accordingly. Update every sheet
# (Overview, zh-cn, de-de) so the file names/dates - and the hyperlink
# text that displays them - reflect this new handoff.

$wb = $excel.ActiveWorkbook

$newGuid = "8775b19e-3c7a-4007-b6a7-3e82fd220f56"
$newHash = "4ce8fb7bdcf3f95cafe2b83f7e6eb133b158738d"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$overviewDate = "2016-48-17 20:48:27"
$zhDate       = "2016-03-17 20:48:21"
$deDate       = "2016-03-17 20:48:27"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = $newMdName
$overview.Range("D2").Value = $overviewDate
$overview.Hyperlinks.Item(1).TextToDisplay = $newMdName

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = $newMdName
$zhcn.Range("D2").Value = $newZhXlfName
$zhcn.Range("E2").Value = $zhDate
$zhcn.Hyperlinks.Item(1).TextToDisplay = $newMdName
$zhcn.Hyperlinks.Item(3).TextToDisplay = $newZhXlfName

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = $newMdName
$dede.Range("D2").Value = $newDeXlfName
$dede.Range("E2").Value = $deDate
$dede.Hyperlinks.Item(1).TextToDisplay = $newMdName
$dede.Hyperlinks.Item(3).TextToDisplay = $newDeXlfName
